$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = 2528.0943407812315
$ws.Range("B6").Value = 1703.9402562987252
$ws.Range("C6").Value = 2002.6567533848233

$ws.Range("A13").Value = 2573.4628886183677
$ws.Range("B13").Value = 2070.9124728367906
$ws.Range("C13").Value = 2069.8969628752407

$ws.Range("A15").Value = 2574.2940135886543
$ws.Range("B15").Value = 1946.9692383072338
$ws.Range("C15").Value = 2186.1572205028783
